# Scheduled runner update: refresh cached Universalis market-data columns
# (H:currentAveragePrice, I/J:.. NQ/HQ, K/L:LevePrice NQ/HQ, M/N:LeveProfit NQ/HQ)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2954.7778
$ws.Range("I9").Value = 3238.2
$ws.Range("K9").Value = 3238.2
$ws.Range("M9").Value = -3069.2
$ws.Range("H40").Value = 6594.4443
$ws.Range("J40").Value = 6850
$ws.Range("L40").Value = 6850
$ws.Range("N40").Value = -7200
$ws.Range("H53").Value = 1853.3077
$ws.Range("J53").Value = 2806.75
$ws.Range("L53").Value = 2806.75
$ws.Range("N53").Value = -4080.75
$ws.Range("H103").Value = 1108.697
$ws.Range("I103").Value = 319.35
$ws.Range("J103").Value = 2323.077
$ws.Range("K103").Value = 958.0500000000001
$ws.Range("L103").Value = 6969.231000000001
$ws.Range("M103").Value = -372.0500000000001
$ws.Range("N103").Value = -8141.231000000001
$ws.Range("H129").Value = 1664.1
$ws.Range("I129").Value = 945.38464
$ws.Range("J129").Value = 2998.8572
$ws.Range("K129").Value = 2836.15392
$ws.Range("L129").Value = 8996.571599999999
$ws.Range("M129").Value = 2163.84608
$ws.Range("N129").Value = -18996.5716
$ws.Range("H138").Value = 2591.5789
$ws.Range("I138").Value = 1086.4849
$ws.Range("J138").Value = 4661.0835
$ws.Range("K138").Value = 3259.4547
$ws.Range("L138").Value = 13983.2505
$ws.Range("M138").Value = 1880.5453
$ws.Range("N138").Value = -24263.2505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 5464.6665
$ws.Range("I28").Value = 5464.6665
$ws.Range("K28").Value = 5464.6665
$ws.Range("M28").Value = -5272.6665
$ws.Range("H32").Value = 4579.5317
$ws.Range("I32").Value = 3838.9048
$ws.Range("J32").Value = 10800.8
$ws.Range("K32").Value = 3838.9048
$ws.Range("L32").Value = 10800.8
$ws.Range("M32").Value = -3551.9048
$ws.Range("N32").Value = -11374.8
$ws.Range("H61").Value = 6256.5
$ws.Range("I61").Value = 4354.4443
$ws.Range("K61").Value = 4354.4443
$ws.Range("M61").Value = -4142.4443
$ws.Range("H99").Value = 5464.6665
$ws.Range("I99").Value = 5464.6665
$ws.Range("K99").Value = 5464.6665
$ws.Range("M99").Value = -2469.6665
$ws.Range("H110").Value = 2985.8
$ws.Range("I110").Value = 2437.4583
$ws.Range("K110").Value = 2437.4583
$ws.Range("M110").Value = -392.4582999999998
$ws.Range("H136").Value = 6256.5
$ws.Range("I136").Value = 4354.4443
$ws.Range("K136").Value = 13063.3329
$ws.Range("M136").Value = -10513.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2001
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -827
$ws.Range("H57").Value = 70780
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 70780
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 70780
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -72220
$ws.Range("H86").Value = 6817.7
$ws.Range("J86").Value = 11898.5
$ws.Range("L86").Value = 11898.5
$ws.Range("N86").Value = -14144.5
$ws.Range("H89").Value = 6817.7
$ws.Range("J89").Value = 11898.5
$ws.Range("L89").Value = 59492.5
$ws.Range("N89").Value = -70724.5
$ws.Range("H136").Value = 70780
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 70780
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 70780
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -80980

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 36510.293
$ws.Range("J31").Value = 82564.5
$ws.Range("L31").Value = 82564.5
$ws.Range("N31").Value = -83154.5
$ws.Range("H34").Value = 36510.293
$ws.Range("J34").Value = 82564.5
$ws.Range("L34").Value = 82564.5
$ws.Range("N34").Value = -82968.5
$ws.Range("H94").Value = 2826.3333
$ws.Range("J94").Value = 3423.5715
$ws.Range("L94").Value = 3423.5715
$ws.Range("N94").Value = -4325.5715
$ws.Range("H102").Value = 40500
$ws.Range("J102").Value = 40666.668
$ws.Range("L102").Value = 40666.668
$ws.Range("N102").Value = -45534.668
$ws.Range("H132").Value = 3150.4443
$ws.Range("I132").Value = 2229.4211
$ws.Range("K132").Value = 6688.263300000001
$ws.Range("M132").Value = -4158.263300000001
$ws.Range("H134").Value = 4132.778
$ws.Range("I134").Value = 2625.1333
$ws.Range("K134").Value = 7875.3999
$ws.Range("M134").Value = -5340.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 85243.414
$ws.Range("I11").Value = 127327.75
$ws.Range("K11").Value = 381983.25
$ws.Range("M11").Value = -381843.25
$ws.Range("H26").Value = 1129.2
$ws.Range("I26").Value = 1013.3333
$ws.Range("K26").Value = 3039.9999
$ws.Range("M26").Value = -2751.9999
$ws.Range("H40").Value = 209.88889
$ws.Range("I40").Value = 41
$ws.Range("J40").Value = 345
$ws.Range("K40").Value = 164
$ws.Range("L40").Value = 1380
$ws.Range("M40").Value = -95
$ws.Range("N40").Value = -1518
$ws.Range("H55").Value = 1880.5385
$ws.Range("I55").Value = 1484.8
$ws.Range("J55").Value = 3199.6667
$ws.Range("K55").Value = 4454.4
$ws.Range("L55").Value = 9599.000100000001
$ws.Range("M55").Value = -4277.4
$ws.Range("N55").Value = -9953.000100000001
$ws.Range("H118").Value = 1557.5
$ws.Range("I118").Value = 481.66666
$ws.Range("J118").Value = 2633.3333
$ws.Range("K118").Value = 1444.99998
$ws.Range("L118").Value = 7899.999899999999
$ws.Range("M118").Value = -201.9999800000001
$ws.Range("N118").Value = -10385.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 847.6923
$ws.Range("I2").Value = 151.66667
$ws.Range("K2").Value = 151.66667
$ws.Range("M2").Value = -38.66667000000001
$ws.Range("H80").Value = 298036.2
$ws.Range("I80").Value = 626906.5
$ws.Range("J80").Value = 5707
$ws.Range("K80").Value = 626906.5
$ws.Range("L80").Value = 5707
$ws.Range("M80").Value = -625908.5
$ws.Range("N80").Value = -7703
$ws.Range("H83").Value = 298036.2
$ws.Range("I83").Value = 626906.5
$ws.Range("J83").Value = 5707
$ws.Range("K83").Value = 3134532.5
$ws.Range("L83").Value = 28535
$ws.Range("M83").Value = -3129540.5
$ws.Range("N83").Value = -38519
$ws.Range("H136").Value = 56356.43
$ws.Range("I136").Value = 57000
$ws.Range("J136").Value = 56306.92
$ws.Range("K136").Value = 171000
$ws.Range("L136").Value = 168920.76
$ws.Range("M136").Value = -168450
$ws.Range("N136").Value = -174020.76

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6571
$ws.Range("J16").Value = 25699
$ws.Range("L16").Value = 25699
$ws.Range("N16").Value = -26039
$ws.Range("H22").Value = 2778.3333
$ws.Range("I22").Value = 1087.1765
$ws.Range("K22").Value = 1087.1765
$ws.Range("M22").Value = -792.1765
$ws.Range("H27").Value = 2778.3333
$ws.Range("I27").Value = 1087.1765
$ws.Range("K27").Value = 1087.1765
$ws.Range("M27").Value = -980.1765
$ws.Range("H46").Value = 4128.375
$ws.Range("J46").Value = 4296.1665
$ws.Range("L46").Value = 4296.1665
$ws.Range("N46").Value = -4672.1665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 170
$ws.Range("I51").Value = 170
$ws.Range("K51").Value = 170
$ws.Range("M51").Value = 340
